# Generate Report for Handback
#
# This script updates the zh-cn / de-de localization-status report after a
# handback: the "Ready for handoff" status becomes "Handed back: in sync
# with en-US", the per-language sheets gain populated "Latest Target File"
# (E) / "Latest Handback File" (F) hyperlink columns, and the "Latest
# Handback DateTime" (G) timestamps are stamped with the handback time.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

# Hyperlink font color/underline used throughout this workbook for link
# cells (matches the existing custom "HyperLink" look: underline + FF6495ED).
$linkColor = 15570276

function Style-AsLink($range) {
    $range.Font.Underline = $true
    $range.Font.Color = $linkColor
}

# ---------------------------------------------------------------------
# Overview sheet: just the status text refresh (B2:C3).
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B2").Value = $newStatus
$wsZh.Range("B3").Value = $newStatus

$zhTargetMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/6d1fbcf9f71fa858e2bc97e097a48cf9225ec2ad/e2e/e2972f55-10dd-48fb-83b2-7caf48ffc041.md"
$zhHandbackXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/345c864cffc74051065243c28b150d061c4f456b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/e2972f55-10dd-48fb-83b2-7caf48ffc041.dfdd3a3e8e7b5877cb446431d5ae0afcaa53cfc9.zh-cn.xlf"
$zhMdName = "e2972f55-10dd-48fb-83b2-7caf48ffc041.md"
$zhXlfName = "e2972f55-10dd-48fb-83b2-7caf48ffc041.dfdd3a3e8e7b5877cb446431d5ae0afcaa53cfc9.zh-cn.xlf"

$wsZh.Range("E2").Value = $zhMdName
$wsZh.Hyperlinks.Add($wsZh.Range("E2"), $zhTargetMdUrl, "", "", $zhMdName)
Style-AsLink $wsZh.Range("E2")

$wsZh.Range("F2").Value = $zhXlfName
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $zhHandbackXlfUrl, "", "", $zhXlfName)
Style-AsLink $wsZh.Range("F2")

$wsZh.Range("E3").Value = $zhMdName
$wsZh.Hyperlinks.Add($wsZh.Range("E3"), $zhTargetMdUrl, "", "", $zhMdName)
Style-AsLink $wsZh.Range("E3")

$wsZh.Range("F3").Value = $zhXlfName
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $zhHandbackXlfUrl, "", "", $zhXlfName)
Style-AsLink $wsZh.Range("F3")

$wsZh.Range("G2").Value = "2016-03-09 13:21:14"
$wsZh.Range("G3").Value = "2016-03-09 13:21:14"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B2").Value = $newStatus
$wsDe.Range("B3").Value = $newStatus

$deTargetMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/6d1fbcf9f71fa858e2bc97e097a48cf9225ec2ad/e2e/e2972f55-10dd-48fb-83b2-7caf48ffc041.md"
$deHandbackXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/140aa2980a6ea81a9f780e6270e25b8cb8fa9774/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/e2972f55-10dd-48fb-83b2-7caf48ffc041.dfdd3a3e8e7b5877cb446431d5ae0afcaa53cfc9.de-de.xlf"
$deMdName = "e2972f55-10dd-48fb-83b2-7caf48ffc041.md"
$deXlfName = "e2972f55-10dd-48fb-83b2-7caf48ffc041.dfdd3a3e8e7b5877cb446431d5ae0afcaa53cfc9.de-de.xlf"

$wsDe.Range("E2").Value = $deMdName
$wsDe.Hyperlinks.Add($wsDe.Range("E2"), $deTargetMdUrl, "", "", $deMdName)
Style-AsLink $wsDe.Range("E2")

$wsDe.Range("F2").Value = $deXlfName
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $deHandbackXlfUrl, "", "", $deXlfName)
Style-AsLink $wsDe.Range("F2")

$wsDe.Range("E3").Value = $deMdName
$wsDe.Hyperlinks.Add($wsDe.Range("E3"), $deTargetMdUrl, "", "", $deMdName)
Style-AsLink $wsDe.Range("E3")

$wsDe.Range("F3").Value = $deXlfName
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $deHandbackXlfUrl, "", "", $deXlfName)
Style-AsLink $wsDe.Range("F3")

$wsDe.Range("G2").Value = "2016-03-09 13:21:35"
$wsDe.Range("G3").Value = "2016-03-09 13:21:35"
